$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "INATIVO - 55.7 meses sem comprar"
$ws.Range("J8").Value = "INATIVO - 18.4 meses sem comprar"
$ws.Range("J18").Value = "INATIVO - 37.8 meses sem comprar"
$ws.Range("J20").Value = "INATIVO - 15.7 meses sem comprar"
$ws.Range("J27").Value = "INATIVO - 19.9 meses sem comprar"
$ws.Range("J29").Value = "INATIVO - 18.0 meses sem comprar"
$ws.Range("J30").Value = "INATIVO - 6.4 meses sem comprar"
$ws.Range("J32").Value = "INATIVO - 22.9 meses sem comprar"
$ws.Range("J35").Value = "INATIVO - 15.0 meses sem comprar"
$ws.Range("J42").Value = "INATIVO - 12.9 meses sem comprar"
$ws.Range("J48").Value = "INATIVO - 2.9 meses sem comprar"
$ws.Range("J53").Value = "INATIVO - 6.5 meses sem comprar"
$ws.Range("J60").Value = "INATIVO - 12.3 meses sem comprar"
$ws.Range("E72").Value = 15
$ws.Range("H72").Value = 45847.9480787037
$ws.Range("I72").Value = 45878.9480787037
$ws.Range("J76").Value = "INATIVO - 33.8 meses sem comprar"
$ws.Range("J85").Value = "INATIVO - 26.6 meses sem comprar"
$ws.Range("J86").Value = "INATIVO - 6.3 meses sem comprar"
$ws.Range("J90").Value = "INATIVO - 5.3 meses sem comprar"
$ws.Range("J94").Value = "INATIVO - 33.8 meses sem comprar"
$ws.Range("J95").Value = "INATIVO - 14.1 meses sem comprar"
$ws.Range("J96").Value = "INATIVO - 19.1 meses sem comprar"
$ws.Range("J97").Value = "INATIVO - 16.8 meses sem comprar"
$ws.Range("J101").Value = "INATIVO - 2.6 meses sem comprar"
$ws.Range("J106").Value = "INATIVO - 25.3 meses sem comprar"
$ws.Range("J107").Value = "INATIVO - 10.9 meses sem comprar"
$ws.Range("J108").Value = "INATIVO - 25.7 meses sem comprar"
$ws.Range("J109").Value = "INATIVO - 16.1 meses sem comprar"
$ws.Range("J110").Value = "INATIVO - 7.0 meses sem comprar"
$ws.Range("J114").Value = "INATIVO - 8.9 meses sem comprar"
$ws.Range("E115").Value = 16473
$ws.Range("H115").Value = 45847.7393287037
$ws.Range("I115").Value = 45848.7393287037
